# Apply the updates described by the diff to the grouped_r_metrics workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "By_Odds_Bin": update rows 2-4 (columns B:F)
# ---------------------------------------------------------------
$wsOdds = $wb.Worksheets.Item("By_Odds_Bin")

$wsOdds.Range("B2").Value = 46
$wsOdds.Range("C2").Value = 4.030000000000001
$wsOdds.Range("D2").Value = 39.03
$wsOdds.Range("E2").Value = -35
$wsOdds.Range("F2").Value = 23.9

$wsOdds.Range("B3").Value = 42
$wsOdds.Range("C3").Value = 23
$wsOdds.Range("D3").Value = 56
$wsOdds.Range("E3").Value = -33
$wsOdds.Range("F3").Value = 21.4

$wsOdds.Range("B4").Value = 3
$wsOdds.Range("C4").Value = -3
$wsOdds.Range("D4").Value = 0
$wsOdds.Range("E4").Value = -3
$wsOdds.Range("F4").Value = 0

# ---------------------------------------------------------------
# Sheet "By_Field_Size": update rows 2-4, relabel rows 7-8, insert a
# new row 9 ("11-13" bin) pushing the old row 9 ("14+") down to row
# 10.
# ---------------------------------------------------------------
$wsField = $wb.Worksheets.Item("By_Field_Size")

$wsField.Range("B2").Value = 26
$wsField.Range("C2").Value = 12.03
$wsField.Range("D2").Value = 31.03
$wsField.Range("E2").Value = -19
$wsField.Range("F2").Value = 26.9

$wsField.Range("B3").Value = 28
$wsField.Range("C3").Value = 2
$wsField.Range("D3").Value = 25
$wsField.Range("E3").Value = -23
$wsField.Range("F3").Value = 17.9

$wsField.Range("B4").Value = 37
$wsField.Range("C4").Value = 10
$wsField.Range("D4").Value = 39
$wsField.Range("E4").Value = -29
$wsField.Range("F4").Value = 21.6

# Relabel the "9-10" / "11-13" bins to the split "9" / "10" bins.
# Go through a TEXT() formula + paste-values so the digit-only labels
# stay text (matching the original inline-string cells) instead of
# being auto-coerced to numbers, without disturbing the cell style.
$wsField.Range("A7").Formula = "=TEXT(9,""0"")"
$wsField.Range("A7").Copy()
$wsField.Range("A7").PasteSpecial(-4163)

$wsField.Range("A8").Formula = "=TEXT(10,""0"")"
$wsField.Range("A8").Copy()
$wsField.Range("A8").PasteSpecial(-4163)

# Insert a fresh row above the old row 9 ("14+"), pushing it to row
# 10, then seed the newly freed row 9 with the "11-13" bin (formatted
# like the sibling zero-bins above it).
$wsField.Rows.Item(9).Insert()
$wsField.Range("A8:F8").Copy($wsField.Range("A9:F9"))

$wsField.Range("A9").Value = "11–13"
$wsField.Range("B9").Value = 0
$wsField.Range("C9").Value = 0
$wsField.Range("D9").Value = 0
$wsField.Range("E9").Value = 0

# ---------------------------------------------------------------
# Sheet "By_Track": replace the single track row with CATTERICK data
# ---------------------------------------------------------------
$wsTrack = $wb.Worksheets.Item("By_Track")

$wsTrack.Range("A2").Value = "CATTERICK"
$wsTrack.Range("B2").Value = 91
$wsTrack.Range("C2").Value = 24.03
$wsTrack.Range("D2").Value = 95.03
$wsTrack.Range("E2").Value = -71
$wsTrack.Range("F2").Value = 22
